# Ultimo cambio despues de publicarlo en el dominio
#
# The sheet lists collection rows grouped by due-date (column L). Several
# pairs of rows that share the same due-date had their data (loan id / A,
# advisor name / B, amount / K) swapped relative to each other, while the
# shared due-date in L stayed put. Re-create that by swapping the A/B/K
# triple between each affected row pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($sheet, $row1, $row2)

    $a1 = $sheet.Range("A$row1").Value2
    $b1 = $sheet.Range("B$row1").Value2
    $k1 = $sheet.Range("K$row1").Value2

    $a2 = $sheet.Range("A$row2").Value2
    $b2 = $sheet.Range("B$row2").Value2
    $k2 = $sheet.Range("K$row2").Value2

    $sheet.Range("A$row1").Value2 = $a2
    $sheet.Range("B$row1").Value2 = $b2
    $sheet.Range("K$row1").Value2 = $k2

    $sheet.Range("A$row2").Value2 = $a1
    $sheet.Range("B$row2").Value2 = $b1
    $sheet.Range("K$row2").Value2 = $k1
}

Swap-RowData $ws 23 24
Swap-RowData $ws 43 44
Swap-RowData $ws 54 55
Swap-RowData $ws 64 65
Swap-RowData $ws 83 84
Swap-RowData $ws 85 86
Swap-RowData $ws 90 91
Swap-RowData $ws 93 94
Swap-RowData $ws 95 96
Swap-RowData $ws 102 104

# Report re-generation timestamp embedded in the sheet's left header.
$ws.PageSetup.LeftHeader = "2024-02-16 06:34:08"
